$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column L (12) to make room for
# "weight_below" and "weight_above", shifting v_alpha/mrpctile/mrdist right.
$ws.Range("L1:M1").EntireColumn.Insert()

# Match style/formatting of neighboring header cell (K1) by copying its
# formatting onto the two freshly-inserted header cells.
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New header labels
$ws.Range("L1").Value = "weight_below"
$ws.Range("M1").Value = "weight_above"

# New data values
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 0.8
$ws.Range("L3").Value = 0.2
$ws.Range("M3").Value = 0.8

# Resize the new columns to fit their (now longer) header text, same as
# Excel's own "best fit" behavior when columns are inserted next to
# auto-fit columns.
$ws.Columns("L").ColumnWidth = 11.850260416666666
$ws.Columns("M").ColumnWidth = 11.744791666666666

# Update selection to match target state
$ws.Range("L3").Select()
